$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (the oldest data point) - this shifts rows 3-36 up to 2-35,
# which naturally gives the correct new values for columns A and B.
$ws.Rows("2:2").Delete()

# Recompute Trad_Prediction (C) as the naive "previous close" forecast,
# and set the new AI_Prediction (D) values for each row.
$ws.Cells.Item(2, 3).Value = 477.8999938964844
$ws.Cells.Item(2, 4).Value = 481.1128510044615
$ws.Cells.Item(3, 3).Value = 477.1199951171875
$ws.Cells.Item(3, 4).Value = 489.8585839254337
$ws.Cells.Item(4, 3).Value = 476.1300048828125
$ws.Cells.Item(4, 4).Value = 473.0241796832669
$ws.Cells.Item(5, 3).Value = 477.3799133300781
$ws.Cells.Item(5, 4).Value = 484.2510013603114
$ws.Cells.Item(6, 3).Value = 477.4675903320312
$ws.Cells.Item(6, 4).Value = 484.5622380781906
$ws.Cells.Item(7, 3).Value = 478.2366943359375
$ws.Cells.Item(7, 4).Value = 472.8172049597523
$ws.Cells.Item(8, 3).Value = 479.2000122070312
$ws.Cells.Item(8, 4).Value = 486.3686712203306
$ws.Cells.Item(9, 3).Value = 477.8099975585938
$ws.Cells.Item(9, 4).Value = 475.5549507719016
$ws.Cells.Item(10, 3).Value = 477.1400146484375
$ws.Cells.Item(10, 4).Value = 498.4437122265703
$ws.Cells.Item(11, 3).Value = 476.5199890136719
$ws.Cells.Item(11, 4).Value = 462.5624233236657
$ws.Cells.Item(12, 3).Value = 476.1199951171875
$ws.Cells.Item(12, 4).Value = 462.6102146511383
$ws.Cells.Item(13, 3).Value = 476.510009765625
$ws.Cells.Item(13, 4).Value = 474.5905620851281
$ws.Cells.Item(14, 3).Value = 476.8599853515625
$ws.Cells.Item(14, 4).Value = 468.8700996682331
$ws.Cells.Item(15, 3).Value = 477.9100036621094
$ws.Cells.Item(15, 4).Value = 491.8447525936203
$ws.Cells.Item(16, 3).Value = 471.960693359375
$ws.Cells.Item(16, 4).Value = 480.7287379051372
$ws.Cells.Item(17, 3).Value = 473.8099975585938
$ws.Cells.Item(17, 4).Value = 454.9223129983978
$ws.Cells.Item(18, 3).Value = 474.8099975585938
$ws.Cells.Item(18, 4).Value = 467.7377256227846
$ws.Cells.Item(19, 3).Value = 475.2000122070312
$ws.Cells.Item(19, 4).Value = 465.0498156032133
$ws.Cells.Item(20, 3).Value = 475.8900146484375
$ws.Cells.Item(20, 4).Value = 472.5471853060204
$ws.Cells.Item(21, 3).Value = 476.3299865722656
$ws.Cells.Item(21, 4).Value = 476.6251869852231
$ws.Cells.Item(22, 3).Value = 475.5599975585938
$ws.Cells.Item(22, 4).Value = 479.8403996255669
$ws.Cells.Item(23, 3).Value = 476.4299926757812
$ws.Cells.Item(23, 4).Value = 476.8463401757896
$ws.Cells.Item(24, 3).Value = 476.9549865722656
$ws.Cells.Item(24, 4).Value = 498.5782136310185
$ws.Cells.Item(25, 3).Value = 477.1176147460938
$ws.Cells.Item(25, 4).Value = 471.685873830732
$ws.Cells.Item(26, 3).Value = 476.6740112304688
$ws.Cells.Item(26, 4).Value = 466.54429113052
$ws.Cells.Item(27, 3).Value = 476.6849975585938
$ws.Cells.Item(27, 4).Value = 475.1415446680317
$ws.Cells.Item(28, 3).Value = 476.9500122070312
$ws.Cells.Item(28, 4).Value = 490.5508740666817
$ws.Cells.Item(29, 3).Value = 476.4100036621094
$ws.Cells.Item(29, 4).Value = 465.0362247423653
$ws.Cells.Item(30, 3).Value = 472.2214965820312
$ws.Cells.Item(30, 4).Value = 491.8666971650148
$ws.Cells.Item(31, 3).Value = 472.5400085449219
$ws.Cells.Item(31, 4).Value = 478.5859702044506
$ws.Cells.Item(32, 3).Value = 471.989990234375
$ws.Cells.Item(32, 4).Value = 473.2872439182959
$ws.Cells.Item(33, 3).Value = 471.8599853515625
$ws.Cells.Item(33, 4).Value = 474.2687753934814
$ws.Cells.Item(34, 3).Value = 470.9312133789062
$ws.Cells.Item(34, 4).Value = 472.2577908070573
$ws.Cells.Item(35, 3).Value = 471.6499938964844
$ws.Cells.Item(35, 4).Value = 474.3269143657882
